$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "NombredelArtista"
$ws.Range("A1").Select() | Out-Null
